$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove excess columns in the DTR summary (rows 15 & 18)
# ---------------------------------------------------------------------------
$ws.Range("I15").ClearContents()

$ws.Range("F18").Value = 3.0
$ws.Range("I18").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add legends section to the per employee report (rows 24-30)
# ---------------------------------------------------------------------------

# "Legends:" header - reuse the big bold/underline title style (same as the
# iRipple, Inc. banner at the top of the sheet).
$ws.Range("E24:P24").Merge()
$ws.Range("E24").Value = "Legends:"
$ws.Range("E24:P24").Font.Name = "Arial"
$ws.Range("E24:P24").Font.Size = 15
$ws.Range("E24:P24").Font.Bold = $true
$ws.Range("E24:P24").Font.Underline = $true

# Legend 1 - blue swatch + request/remark note
$ws.Range("E25:E26").Merge()
$ws.Range("E25:E26").Interior.Color = 0xCCA329

$ws.Range("F25:P26").Merge()
$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."
$ws.Range("F25:P26").Font.Bold = $true
$ws.Range("F25:P26").Font.Underline = $true

# Legend 2 - orange swatch + half-day note
$ws.Range("E27:E28").Merge()
$ws.Range("E27:E28").Interior.Color = 0x66CCFF

$ws.Range("F27:P28").Merge()
$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."
$ws.Range("F27:P28").Font.Bold = $true
$ws.Range("F27:P28").Font.Underline = $true

# Legend 3 - red swatch + absent note
$ws.Range("E29:E30").Merge()
$ws.Range("E29:E30").Interior.Color = 0x5E5EDF

$ws.Range("F29:P30").Merge()
$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
$ws.Range("F29:P30").Font.Bold = $true
$ws.Range("F29:P30").Font.Underline = $true
